# Regenerate save_data K column (Strike# -> K), updating the G column
# values for rows 2-17 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 5
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 2
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
